$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (D, E, F) before the existing "Terms Typically
# Offered" column, shifting it from D to G.
$ws.Range("D1:F1").EntireColumn.Insert()

# Header row
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Fill the new columns with "NA" for every data row (rows 2-14)
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 4).Value = "NA"
    $ws.Cells.Item($row, 5).Value = "NA"
    $ws.Cells.Item($row, 6).Value = "NA"
}
